$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.904.64"
$ws.Range("D3").Value = "'1.635.81"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'211.70"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'23.44"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "'1.867.99"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "'1.647.04"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'65.39"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'27.919.82"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'228.80"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'0.0₃0720"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "'4.36"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'10.08"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").Value = "'155.52"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D30").Value = "'1.18"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "'1.396.63"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "'0.560"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").Value = "'66.04"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "'1.777.28"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "  +1.58%  "
